$wb = $excel.ActiveWorkbook

# --- DBS sheet: append 4 new lookup rows (AgStatusCode variants of the
#     existing AgCurInd rows) right after the current last row (11) ---
$ws2 = $wb.Worksheets.Item("DBS")

# Column A first (FunNm), then column B (讀取Key條件), then column C
# (其他ORDER條件) - matches the shared-string insertion order of the
# original edit (new strings 328-331 are the A-column names, 332-335 are
# the B-column conditions).
$ws2.Range("A12").Value = "findCenterCodeAndAgStatusCode"
$ws2.Range("A13").Value = "EmployeeNoLikeAndAgStatusCode"
$ws2.Range("A14").Value = "findFullnameLikeAndAgStatusCode"
$ws2.Range("A15").Value = "findEmployeeNoAndAgStatusCode"

$ws2.Range("B12").Value = "CenterCode = , AND AgStatusCode = "
$ws2.Range("B13").Value = "EmployeeNo % , AND AgStatusCode = "
$ws2.Range("B14").Value = "Fullname % , AND AgStatusCode = "
$ws2.Range("B15").Value = "EmployeeNo >= ,AND EmployeeNo <= , AND AgStatusCode = "

$ws2.Range("C12").Value = "EmployeeNo Asc"
$ws2.Range("C13").Value = "EmployeeNo Asc"
$ws2.Range("C14").Value = "EmployeeNo Asc"
$ws2.Range("C15").Value = "EmployeeNo ASC"

# --- view state: DBS becomes the active/visible tab, scrolled down so the
#     newly added rows are in view, with B19 as the active cell in the
#     frozen bottom pane ---
$ws1 = $wb.Worksheets.Item("DBD")
$ws1.Activate() | Out-Null
$ws1.Range("B35").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B19").Select() | Out-Null
